$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the dates in column A for rows 9-15 forward by one week (7 days)
$ws.Range("A9").Value = 45684
$ws.Range("A10").Value = 45691
$ws.Range("A11").Value = 45698
$ws.Range("A12").Value = 45705
$ws.Range("A13").Value = 45712
$ws.Range("A14").Value = 45719
$ws.Range("A15").Value = 45726

# Fill in values for row 9 (columns G, H, I, L, M, O, P, Q)
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 5
$ws.Range("I9").Value = 5
$ws.Range("L9").Value = 5
$ws.Range("M9").Value = 5
$ws.Range("O9").Value = 5
$ws.Range("P9").Value = 5
$ws.Range("Q9").Value = 5

# Update the selection to match the diff
$ws.Range("O9").Select()
